# Generate Report for Handback
# A new handback run completed for file "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md",
# refreshing the handoff/handback timestamps recorded for the zh-cn and de-de
# targets (and the rolled-up "Latest HO Xliff Generate Date" on the Overview sheet).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for 1178f3c0-...md (row 2)
$wsOverview.Range("G2").Value = "2016-08-31 07:44:02"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for 1178f3c0-...md (row 2)
$wsZhCn.Range("H2").Value = "2016-08-31 07:43:50"
$wsZhCn.Range("K2").Value = "2016-08-31 07:44:34"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime for 1178f3c0-...md (row 2)
$wsDeDe.Range("H2").Value = "2016-08-31 07:44:02"
$wsDeDe.Range("K2").Value = "2016-08-31 07:44:52"
